$wb = $excel.ActiveWorkbook

# --- "survey" sheet: remove the plot_id question row ---
# The form no longer asks for a manually entered plot id (instance_name is
# used instead), so the "integer / plot_id / Enter the id of plot:" row is
# deleted and everything below it shifts up.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Delete()
$ws1.Range("B7").Select() | Out-Null

# --- "settings" sheet: add table_id + disableSwipeNavigation settings ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A6").Value = "table_id"
$ws3.Range("B6").Value = "plot"

$ws3.Range("A7").Value = "disableSwipeNavigation"
$ws3.Range("A7").Font.Color = 0
$ws3.Range("A7").WrapText = $true
$ws3.Range("A7").Font.Name = "Times New Roman"

$ws3.Range("B7").Value = $true
$ws3.Range("B7").Font.Color = 0
$ws3.Range("B7").Font.Size = 10
$ws3.Range("B7").WrapText = $true
$ws3.Range("B7").Font.Name = "Arial"

$ws3.Rows.Item(7).RowHeight = 30

$ws3.Range("A7:B7").Select() | Out-Null
